# Updated symbol list on Wed Dec 21 16:11:04 UTC 2022 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Hora" (column G) columns for the
# crypto symbol list. Values are written with a leading apostrophe so the
# numeric-looking price strings stay stored as text (matching the sheet's
# existing inline-string layout) instead of being coerced to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Price" values (column D) per row. Rows not listed here keep their
# current price untouched (e.g. rows whose price is "--" or unchanged).
$priceUpdates = @{
    2  = "248.73"
    3  = "22.63"
    4  = "5.272"
    5  = "0.05701"
    6  = "3.406"
    7  = "6.340"
    8  = "0.8090"
    9  = "0.8941"
    10 = "0.1423"
    11 = "0.07443"
    12 = "0.03053"
    13 = "0.03099"
    14 = "0.09409"
    15 = "3.868"
    16 = "0.001571"
    17 = "0.04788"
    18 = "0.01827"
    19 = "0.0005800"
    20 = "0.006422"
    21 = "0.004982"
    22 = "0.0009959"
    24 = "3.689"
    25 = "2.163"
    26 = "0.3259"
    27 = "0.1369"
    40 = "0.03974"
    41 = "0.006800"
    42 = "0.1070"
    43 = "0.003200"
    44 = "0.007827"
    45 = "0.00005576"
    48 = "0.1979"
}

foreach ($row in $priceUpdates.Keys) {
    $cell = $ws.Cells.Item($row, 4)
    $cell.Value = "'" + $priceUpdates[$row]
    # Drop the quote-prefix style the apostrophe trick stamps on the cell so
    # formatting stays identical to the original (text value, default style).
    $cell.Style = "Normal"
}

# "Hora" (column G) bumps from 15 to 16 for every data row (2 through 51).
for ($row = 2; $row -le 51; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $cell.Value = "'16"
    $cell.Style = "Normal"
}
